# Daily scrape update - 2025-09-12 03:04:21 UTC
# Refreshes the Global Talent opportunity listing: rows 2-12 get new
# opportunity data and the old rows 13-17 are dropped (sheet shrinks from
# A1:H17 to A1:H12). A few column widths are also tweaked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New opportunity rows (2-12) -------------------------------------------------
$rows = @(
    @{ Id = "1327574"; Title = "Controller and Performance Manager"; Country = "Lisboa, Portugal"; Applicants = "6 applicants"; Duration = "6 - 18 Months"; Org = "Portway, Handling de Portugal" },
    @{ Id = "1327541"; Title = "Software Developer"; Country = "União das freguesias de Cascais e Estoril, Portugal"; Applicants = "9 applicants"; Duration = "3 - 6 Months"; Org = "Dark Cloud" },
    @{ Id = "1327539"; Title = "Graphic Designer"; Country = "União das freguesias de Cascais e Estoril, Portugal"; Applicants = "0 applicants"; Duration = "3 - 6 Months"; Org = "Dark Cloud" },
    @{ Id = "1327527"; Title = "Social Media Content Creator"; Country = "Galle, Sri Lanka"; Applicants = "8 applicants"; Duration = "3 - 6 Months"; Org = "Radisson Collection Resort Galle" },
    @{ Id = "1327518"; Title = "Business Development Intern"; Country = "Malabe, Sri Lanka"; Applicants = "3 applicants"; Duration = "3 - 6 Months"; Org = "ZILLIONe Technologies Private Limited" },
    @{ Id = "1327273"; Title = "Aged Care Nurse"; Country = "Melbourne VIC, Australia"; Applicants = "26 applicants"; Duration = "6 - 18 Months"; Org = "Opulence College Pty Ltd" },
    @{ Id = "1327236"; Title = "Videographer & Video Editor"; Country = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"; Applicants = "0 applicants"; Duration = "3 - 6 Months"; Org = "Karcel" },
    @{ Id = "1327232"; Title = "Content Creator (Storyteller & social media Maven)"; Country = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"; Applicants = "0 applicants"; Duration = "3 - 6 Months"; Org = "Karcel" },
    @{ Id = "1327124"; Title = "Graphic Designer"; Country = "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt"; Applicants = "0 applicants"; Duration = "3 - 6 Months"; Org = "Karcel" },
    @{ Id = "1326658"; Title = "Business Development Intern (Japanese Speaking Individuals Only)"; Country = "Malabe, Sri Lanka"; Applicants = "19 applicants"; Duration = "3 - 6 Months"; Org = "Creative Technology Solutions (Private) Limited" },
    @{ Id = "1326003"; Title = "European Epilepsy Trainee ( ONLY EU)"; Country = "Bruxelles, Belgio"; Applicants = "43 applicants"; Duration = "6 - 18 Months"; Org = "UCB" }
)

# Column A holds opportunity IDs that look numeric but must stay text (as in
# the source data). Format the range as Text before writing so the digit
# strings aren't auto-coerced to numbers, then clear the formatting again
# (the source cells carry no explicit style) while the stored values remain
# text.
$idRange = $ws.Range("A2:A12")
$idRange.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Id
    $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $row.Id
    $ws.Cells.Item($r, 3).Value = $row.Title
    $ws.Cells.Item($r, 4).Value = $row.Country
    $ws.Cells.Item($r, 5).Value = "No"
    $ws.Cells.Item($r, 6).Value = $row.Applicants
    $ws.Cells.Item($r, 7).Value = $row.Duration
    $ws.Cells.Item($r, 8).Value = $row.Org
    $r = $r + 1
}

$idRange.ClearFormats()

# Drop the old rows 13-17 (sheet now ends at row 12) --------------------------
$ws.Range("A13:H17").Delete()

# Column width tweaks -----------------------------------------------------
# (ColumnWidth setter adds Excel's standard ~5/6 character padding, so we
# subtract it here to land exactly on the target stored widths.)
$ws.Columns.Item(3).ColumnWidth = 67 - 5/6
$ws.Columns.Item(4).ColumnWidth = 65 - 5/6
$ws.Columns.Item(6).ColumnWidth = 16 - 5/6
$ws.Columns.Item(8).ColumnWidth = 50 - 5/6
